$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G: all rows 2-51 change from "21" to "22" (Hora column)
$ws.Range("G2:G51").Value = "'22"
$ws.Range("G2:G51").Style = "Normal"

# Column D: individual price updates
$ws.Range("D2").Value = "'286.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'21.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'6.459"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.06379"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.601"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'1.558"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'6.571"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8247"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Value = "'0.1674"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.08702"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03699"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.03209"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.09187"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.707"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.001652"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.04749"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006156"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.006286"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.001073"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Value = "'3.783"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.321"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.3356"
$ws.Range("D25").Style = "Normal"
$ws.Range("D40").Value = "'0.04799"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.007183"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.004508"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.1114"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.01170"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006956"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.9349"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.005474"
$ws.Range("D48").Style = "Normal"
